$d = $word.ActiveDocument

# Locate the paragraph right after "LOM3234: Optica Fisica (Requisito)"
# (the blank paragraph that precedes the "Ver no Jupiter ..." line) and the
# paragraph holding the trailing copyright notice ("(c) 2020 ... Contact:
# luizeleno@usp.br ..."). Everything from the start of the former through
# the end (incl. paragraph mark) of the latter must be removed, collapsing
# the three trailing paragraphs ("" / "Ver no Jupiter ..." / "(c) 2020 ...")
# into nothing while leaving the rest of the document untouched.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "*LOM3234*") {
        # The blank paragraph right after this one is where the deletion begins.
        $startPara = $i + 1
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $endPara = $i
    }
}

if ($startPara -ne $null -and $endPara -ne $null -and $startPara -le $endPara) {
    $rStart = $d.Paragraphs.Item($startPara).Range.Start
    $rEnd = $d.Paragraphs.Item($endPara).Range.End
    $r = $d.Range($rStart, $rEnd)
    $r.Delete()
}
